$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

# Reorder the comma-separated "Recorded By" author list for the rows
# touched by the upstream sync (reverses the display order of the same
# set of authors; single-author cells are untouched).
$updates = @{
    "G2" = "system, System, backup@backdoor.com"
    "G3" = "dnasr281@gmail.com, System"
    "G4" = "System, backup@backdoor.com"
    "G5" = "System, backup@backdoor.com"
    "G6" = "dnasr281@gmail.com, System"
    "G7" = "System, admin@admin.com"
    "G8" = "System, backup@backdoor.com"
    "G10" = "dnasr281@gmail.com, System"
    "G11" = "dnasr281@gmail.com, System"
    "G12" = "dnasr281@gmail.com, System"
    "G13" = "dnasr281@gmail.com, System"
    "G14" = "dnasr281@gmail.com, System"
    "G15" = "dnasr281@gmail.com, System"
    "G17" = "dnasr281@gmail.com, System"
    "G18" = "dnasr281@gmail.com, System"
    "G19" = "dnasr281@gmail.com, System"
    "G20" = "dnasr281@gmail.com, System"
    "G21" = "dnasr281@gmail.com, System"
    "G22" = "dnasr281@gmail.com, System"
    "G24" = "dnasr281@gmail.com, System"
    "G26" = "dnasr281@gmail.com, System"
    "G28" = "system, System, backup@backdoor.com"
    "G29" = "dnasr281@gmail.com, System"
    "G30" = "System, backup@backdoor.com"
    "G31" = "System, backup@backdoor.com"
    "G32" = "dnasr281@gmail.com, System"
    "G33" = "System, admin@admin.com"
    "G34" = "System, backup@backdoor.com"
    "G36" = "dnasr281@gmail.com, System"
    "G37" = "dnasr281@gmail.com, System"
    "G38" = "dnasr281@gmail.com, System"
    "G39" = "dnasr281@gmail.com, System"
    "G40" = "dnasr281@gmail.com, System"
    "G41" = "dnasr281@gmail.com, System"
    "G43" = "dnasr281@gmail.com, System"
    "G44" = "dnasr281@gmail.com, System"
    "G45" = "dnasr281@gmail.com, System"
    "G46" = "dnasr281@gmail.com, System"
    "G47" = "dnasr281@gmail.com, System"
    "G48" = "dnasr281@gmail.com, System"
    "G50" = "dnasr281@gmail.com, System"
    "G52" = "dnasr281@gmail.com, System"
    "G54" = "system, System, backup@backdoor.com"
    "G55" = "dnasr281@gmail.com, System"
    "G56" = "System, backup@backdoor.com"
    "G57" = "System, backup@backdoor.com"
    "G58" = "dnasr281@gmail.com, System"
    "G59" = "System, admin@admin.com"
    "G60" = "System, backup@backdoor.com"
    "G62" = "dnasr281@gmail.com, System"
    "G63" = "dnasr281@gmail.com, System"
    "G64" = "dnasr281@gmail.com, System"
    "G65" = "dnasr281@gmail.com, System"
    "G66" = "dnasr281@gmail.com, System"
    "G67" = "dnasr281@gmail.com, System"
    "G69" = "dnasr281@gmail.com, System"
    "G70" = "dnasr281@gmail.com, System"
    "G71" = "dnasr281@gmail.com, System"
    "G72" = "dnasr281@gmail.com, System"
    "G73" = "dnasr281@gmail.com, System"
    "G74" = "dnasr281@gmail.com, System"
    "G76" = "dnasr281@gmail.com, System"
    "G78" = "dnasr281@gmail.com, System"
    "G80" = "System, backup@backdoor.com"
    "G81" = "System, backup@backdoor.com"
    "G82" = "System, backup@backdoor.com"
    "G83" = "dnasr281@gmail.com, System"
    "G84" = "dnasr281@gmail.com, System"
    "G85" = "dnasr281@gmail.com, System"
    "G86" = "dnasr281@gmail.com, System"
    "G87" = "dnasr281@gmail.com, admin@admin.com"
    "G90" = "dnasr281@gmail.com, System"
    "G92" = "dnasr281@gmail.com, System"
    "G93" = "dnasr281@gmail.com, System"
    "G94" = "dnasr281@gmail.com, System"
    "G96" = "dnasr281@gmail.com, System"
    "G99" = "dnasr281@gmail.com, System"
    "G101" = "dnasr281@gmail.com, System"
    "G106" = "System, backup@backdoor.com"
    "G107" = "System, backup@backdoor.com"
    "G108" = "System, backup@backdoor.com"
    "G109" = "dnasr281@gmail.com, System"
    "G110" = "dnasr281@gmail.com, System"
    "G111" = "dnasr281@gmail.com, System"
    "G112" = "dnasr281@gmail.com, System"
    "G113" = "dnasr281@gmail.com, admin@admin.com"
    "G116" = "dnasr281@gmail.com, System"
    "G118" = "dnasr281@gmail.com, System"
    "G119" = "dnasr281@gmail.com, System"
    "G120" = "dnasr281@gmail.com, System"
    "G122" = "dnasr281@gmail.com, System"
    "G125" = "dnasr281@gmail.com, System"
    "G127" = "dnasr281@gmail.com, System"
    "G132" = "System, backup@backdoor.com"
    "G133" = "System, backup@backdoor.com"
    "G134" = "System, backup@backdoor.com"
    "G135" = "dnasr281@gmail.com, System"
    "G136" = "dnasr281@gmail.com, System"
    "G137" = "dnasr281@gmail.com, System"
    "G138" = "dnasr281@gmail.com, System"
    "G139" = "dnasr281@gmail.com, admin@admin.com"
    "G142" = "dnasr281@gmail.com, System"
    "G144" = "dnasr281@gmail.com, System"
    "G145" = "dnasr281@gmail.com, System"
    "G146" = "dnasr281@gmail.com, System"
    "G148" = "dnasr281@gmail.com, System"
    "G151" = "dnasr281@gmail.com, System"
    "G153" = "dnasr281@gmail.com, System"
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
